# Automatic update of files.
# Rows 3-8 (data rows) have their observation-level content rotated between
# row positions while the header row and row 2 stay untouched. This mirrors
# a re-sort/re-sync of the underlying "Artfynd" export where each record's
# Id, coordinates (Ost/Nord), accuracy (Noggrannhet), species fields and
# reporter names move to a different row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (content becomes what used to be in row 7) ---
$ws.Range("A3").Value = 108367417
$ws.Range("Q3").Value = 663132.0536545257
$ws.Range("R3").Value = 6705596.601038971
$ws.Range("S3").Value = 15
$ws.Range("AC3").Value = ""
$ws.Range("AW3").Value = "Isac Carlsson"
$ws.Range("AX3").Value = "Isac Carlsson"

# --- Row 4 (content becomes what used to be in row 3) ---
$ws.Range("A4").Value = 108369210
$ws.Range("Q4").Value = 663162.7064135609
$ws.Range("R4").Value = 6705981.337152475
$ws.Range("S4").Value = 25
$ws.Range("AC4").Value = "#SAKNAS!"
$ws.Range("AW4").Value = "Nadja Nilsson"
$ws.Range("AX4").Value = "Nadja Nilsson"

# --- Row 5 (content becomes what used to be in row 4) ---
$ws.Range("A5").Value = 108367419
$ws.Range("Q5").Value = 663092.6272863077
$ws.Range("R5").Value = 6705966.322238538
$ws.Range("S5").Value = 15
$ws.Range("AC5").Value = ""
$ws.Range("AW5").Value = "Isac Carlsson"
$ws.Range("AX5").Value = "Isac Carlsson"

# --- Row 6 (content becomes what used to be in row 5) ---
$ws.Range("A6").Value = 108369211
$ws.Range("B6").Value = 98520
$ws.Range("E6").Value = 222498
$ws.Range("F6").Value = "Blåsippa"
$ws.Range("G6").Value = "Hepatica nobilis"
$ws.Range("H6").Value = "Schreb."
$ws.Range("Q6").Value = 663155.2423936725
$ws.Range("R6").Value = 6706004.217932139
$ws.Range("S6").Value = 25
$ws.Range("AC6").Value = "#SAKNAS!"
$ws.Range("AW6").Value = "Nadja Nilsson"
$ws.Range("AX6").Value = "Nadja Nilsson"

# --- Row 7 (content becomes what used to be in row 8) ---
$ws.Range("A7").Value = 108367421
$ws.Range("Q7").Value = 663073.8261583259
$ws.Range("R7").Value = 6705549.019092269

# --- Row 8 (content becomes what used to be in row 6) ---
$ws.Range("A8").Value = 108367415
$ws.Range("B8").Value = 95511
$ws.Range("E8").Value = 221944
$ws.Range("F8").Value = "Lopplummer"
$ws.Range("G8").Value = "Huperzia selago"
$ws.Range("H8").Value = "(L.) Bernh. ex Schrank & Mart."
$ws.Range("Q8").Value = 663125.9516252303
$ws.Range("R8").Value = 6705852.723180643
